$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 78, pushing existing rows 78+ down by 2.
$ws.Rows("78:79").Insert()

# --- New row 78 ---
$ws.Range("A78").Value = 10
$ws.Range("B78").Value = "Vega Modelo de Temuco"
$ws.Range("C78").Value = "La Araucanía"
$ws.Range("D78").Value = "2021-08-05"
$ws.Range("E78").Value = 9
$ws.Range("F78").Value = 100112040
$ws.Range("G78").Value = "Cilantro"
$ws.Range("H78").Value = "Sin especificar"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 20
$ws.Range("K78").Value = 5000
$ws.Range("L78").Value = 5000
$ws.Range("M78").Value = 5000
$ws.Range("N78").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O78").Value = "Provincia de Cautín"
$ws.Range("P78").Value = 2500
$ws.Range("Q78").Value = 2
$ws.Range("R78").Value = "Hortaliza"

# --- New row 79 ---
$ws.Range("A79").Value = 10
$ws.Range("B79").Value = "Vega Modelo de Temuco"
$ws.Range("C79").Value = "La Araucanía"
$ws.Range("D79").Value = "2021-08-05"
$ws.Range("E79").Value = 9
$ws.Range("F79").Value = 100112040
$ws.Range("G79").Value = "Cilantro"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 80
$ws.Range("K79").Value = 4300
$ws.Range("L79").Value = 4300
$ws.Range("M79").Value = 4300
$ws.Range("N79").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O79").Value = "Región Metropolitana"
$ws.Range("P79").Value = 2150
$ws.Range("Q79").Value = 2
$ws.Range("R79").Value = "Hortaliza"
